# Sprint 4 Backlog - Burndown
# Commit: "Updated sprint backlog again to include my time for adding all
# items from shopping list to pantry. Included my time again due to merge
# conflits."
#
# Row 12 corresponds to the task "Compete functionality to add all shopping
# list ingredients to pantry (desktop)" (Assigned Team Member = Janera).
# This records Janera's actual time spent (0.75) and marks the task as
# completed by her, with 0 remaining in Week 1 / Week 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E12").Value = 0.75      # Actual Time
$ws.Range("F12").Value = "Janera"  # Completed By
$ws.Range("H12").Value = 0         # Amount Remaining After Week 1
$ws.Range("I12").Value = 0         # Amount Remaining After Week 2

# Leave the active selection on I12, matching the last-edited cell.
$ws.Range("I12").Select()
